$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @"
Aug,0,0,0,0,0,0,0,7.853740234375,8.87140869140625,11.526796875,23.8244970703125,29.21532958984375,15.1430810546875,10.426328125,6.2696533203125,5.365075683593751,4.96929931640625,4.2556201171875,5.142919921875,2.0255908203125,2.9773828125,5.30373046875,4.8036083984375,3.67880859375,9.74630126953125,12.3238671875,5.17105224609375,2.173828125,0.95654296875,1.4128662109375,4.5983154296875
Sep,8.694785156250001,14.22405029296875,8.594409179687501,1.32636474609375,1.27501708984375,1.08232421875,1.607080078125,4.892268066406251,7.02943359375,5.257075195312501,4.13421630859375,5.833947753906251,7.4504345703125,9.0090087890625,21.4332421875,15.8147314453125,3.54728515625,4.507861328125,5.367509765625,3.07078125,4.550244140625,5.023730468749999,3.3130908203125,1.9299560546875,1.47852783203125,1.48900146484375,1.5380859375,2.59954833984375,5.1651611328125,7.800708007812499,0
Oct,0.2548828125,0.12158203125,0.12158203125,0.0400390625,2.156494140625,13.459716796875,12.439453125,5.57666015625,0.087158203125,0.087158203125,0.087158203125,3.608642578125,14.088623046875,14.088623046875,8.13134765625,0.63916015625,0.93359375,1.880859375,9.7255859375,15.91357421875,15.997314453125,8.23388671875,2.232421875,0.090087890625,1.7314453125,3.621826171875,4.779052734375,4.779052734375,2.0625,0.0400390625,0.16015625
Nov,0.090087890625,0.13623046875,0.31640625,0.49560546875,0.1142578125,0.0400390625,0.0400390625,0.037109375,0.04833984375,0.090087890625,0.090087890625,0.081298828125,0.0341796875,0.008544921875,0,0.0368896484375,0.04453125,0.41455078125,0.303515625,0.18212890625,0.29296875,0.03466796875,0.146064453125,0.09424072265625001,0.168818359375,0.08340087890625,0.0371875,0.009033203125,0.03466796875,0.036474609375,0
Dec,0.09991699218750001,0.0966796875,0.009521484375,0.009033203125,0.0166015625,0.0166015625,0.017431640625,0.064453125,0.064453125,0.0380859375,0.235400390625,0.235400390625,0.1510009765625,0.009521484375,0,0,0,0.015771484375,0.015771484375,0.0095849609375,0.0382080078125,0.0421142578125,0.011474609375,0.009521484375,0.0361328125,0.0361328125,0.009033203125,0,0.011474609375,0.011474609375,0.0369140625
Jan,0.0359375,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0.145751953125,0.145751953125,0.116455078125,0,0.081298828125,0,0,0,0,0,0,0,0
Feb,0,0,0,0,0,0,0,0.00830078125,0.009130859375,0.009130859375,0,0,0,0,0.07275390625,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0
Mar,0,0,0.036376953125,0.031494140625,0.046142578125,0.11865234375,0.151611328125,0.19189453125,0.19189453125,0.136962890625,0.1865234375,0.67138671875,0.67138671875,0.872802734375,1.13525390625,0.07958984375,0.011962890625,0,0,0,0,0,0,0,0,0,0,0,0,0.0126953125,0.010009765625
Apr,0.24169921875,0.24169921875,0.34716796875,0.0510009765625,5.411486816406249,7.273251953125,4.535,3.85234375,2.5365234375,0.2255859375,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0
Jun,0,0,0,0,0,0,0,0,0,0,0,0.010498046875,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0
Jul,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0.08740234375,0.10693359375,0.10693359375,0.08740234375,0,0,0
"@

$rows = $data -split "`n"

$startRow = 2
$lastExistingRow = 9

for ($i = 0; $i -lt $rows.Length; $i++) {
    $line = $rows[$i].Trim()
    if ($line.Length -eq 0) { continue }
    $fields = $line -split ","
    $r = $startRow + $i

    if ($r -gt $lastExistingRow) {
        # New row beyond the original data block - clone month-label formatting
        # (border + bold + centered style) from the last existing labeled row.
        $ws.Cells.Item($lastExistingRow, 1).Copy()
        $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    }

    $ws.Cells.Item($r, 1).Value = $fields[0]
    for ($c = 1; $c -lt $fields.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = [double]$fields[$c]
    }
}

Write-Output "Updated $($rows.Length) rows"
